$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Column V header (date), same style/format as U1
$ws.Cells.Item(1, 21).Copy($ws.Cells.Item(1, 22))
$ws.Cells.Item(1, 22).Value = Get-Date -Year 2013 -Month 6 -Day 4 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Column V rows 2-66: mirror column U, with two exceptions (row 12 -> "?", row 49 -> "OK")
$ws.Cells.Item(2, 22).Value = "OK"
$ws.Cells.Item(3, 22).Value = "OK"
$ws.Cells.Item(4, 22).Value = "OK"
$ws.Cells.Item(5, 22).Value = "OK"
$ws.Cells.Item(6, 22).Value = "OK"
$ws.Cells.Item(7, 22).Value = "OK"
$ws.Cells.Item(8, 22).Value = "OK"
$ws.Cells.Item(9, 22).Value = "OK"
$ws.Cells.Item(10, 22).Value = "OK"
$ws.Cells.Item(11, 22).Value = "OK"
$ws.Cells.Item(12, 22).Value = "?"
$ws.Cells.Item(13, 22).Value = "OK"
$ws.Cells.Item(14, 22).Value = "OK"
$ws.Cells.Item(15, 22).Value = "OK"
$ws.Cells.Item(16, 22).Value = "OK"
$ws.Cells.Item(17, 22).Value = "NG"
$ws.Cells.Item(18, 22).Value = "OK"
$ws.Cells.Item(19, 22).Value = "OK"
$ws.Cells.Item(20, 22).Value = "OK"
$ws.Cells.Item(21, 22).Value = "OK"
$ws.Cells.Item(22, 22).Value = "OK"
$ws.Cells.Item(23, 22).Value = "OK"
$ws.Cells.Item(24, 22).Value = "OK"
$ws.Cells.Item(25, 22).Value = "OK"
$ws.Cells.Item(26, 22).Value = "OK"
$ws.Cells.Item(27, 22).Value = "OK"
$ws.Cells.Item(28, 22).Value = "?"
$ws.Cells.Item(29, 22).Value = "?"
$ws.Cells.Item(30, 22).Value = "OK"
$ws.Cells.Item(31, 22).Value = "OK"
$ws.Cells.Item(32, 22).Value = "OK"
$ws.Cells.Item(33, 22).Value = "OK"
$ws.Cells.Item(34, 22).Value = "OK"
$ws.Cells.Item(35, 22).Value = "OK"
$ws.Cells.Item(36, 22).Value = "OK"
$ws.Cells.Item(37, 22).Value = "OK"
$ws.Cells.Item(38, 22).Value = "OK"
$ws.Cells.Item(39, 22).Value = "OK"
$ws.Cells.Item(40, 22).Value = "?"
$ws.Cells.Item(41, 22).Value = "?"
$ws.Cells.Item(42, 22).Value = "OK"
$ws.Cells.Item(43, 22).Value = "OK"
$ws.Cells.Item(44, 22).Value = "OK"
$ws.Cells.Item(45, 22).Value = "OK"
$ws.Cells.Item(46, 22).Value = "OK"
$ws.Cells.Item(47, 22).Value = "OK"
$ws.Cells.Item(48, 22).Value = "?"
$ws.Cells.Item(49, 22).Value = "OK"
$ws.Cells.Item(50, 22).Value = "OK"
$ws.Cells.Item(51, 22).Value = "OK"
$ws.Cells.Item(52, 22).Value = "OK"
$ws.Cells.Item(53, 22).Value = "OK"
$ws.Cells.Item(54, 22).Value = "OK"
$ws.Cells.Item(55, 22).Value = "OK"
$ws.Cells.Item(56, 22).Value = "OK"
$ws.Cells.Item(57, 22).Value = "?"
$ws.Cells.Item(58, 22).Value = "OK"
$ws.Cells.Item(59, 22).Value = "OK"
$ws.Cells.Item(60, 22).Value = "OK"
$ws.Cells.Item(61, 22).Value = "OK"
$ws.Cells.Item(62, 22).Value = "OK"
$ws.Cells.Item(63, 22).Value = "OK"
$ws.Cells.Item(64, 22).Value = "OK"
$ws.Cells.Item(65, 22).Value = "OK"
$ws.Cells.Item(66, 22).Value = "OK"

# Column width for V to match U, and widen column B
$ws.Columns.Item(22).ColumnWidth = 9.33
$ws.Columns.Item(2).ColumnWidth = 86.85

# Update selection to V5 and reset the view to the top-left
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("V5").Select()
